$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D9").Value = 1.577923291743519
$ws.Range("D10").Value = 1.509905915402546
$ws.Range("D11").Value = 1.449960134423589
$ws.Range("D12").Value = 0.8671894150874061

$ws.Range("D20").Value = 1.457453392513047
$ws.Range("D21").Value = 1.323365514853331
$ws.Range("D22").Value = 1.351078349560582

$ws.Range("E31").Value = 13
$ws.Range("E34").Value = 25
